# prise en compte Mbus liste pour config WIT
# Insert two new detail rows ("Compteur Modbus 1" and "Compteur Impul 1")
# just above the totals row, and refresh the totals row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The totals row is currently row 35. Insert two blank rows above it
# (formats are inherited from the row above, i.e. the "Compteurs" style).
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).Insert()

# New row 35: Compteurs / Compteur Modbus 1
$ws.Range("A35").Value = "Compteurs"
$ws.Range("B35").Value = "Compteur Modbus 1"
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 10

# New row 36: Compteurs / Compteur Impul 1
$ws.Range("A36").Value = "Compteurs"
$ws.Range("B36").Value = "Compteur Impul 1"
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0

# Totals row, now shifted down to row 37: refresh label and totals
$ws.Range("A37").Value = "zzzzeTOTAL"
$ws.Range("B37").Value = "TOTAUX (53 points)"
$ws.Range("C37").Value = 9
$ws.Range("D37").Value = 12
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 10
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 10
